$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.07664366666667
$ws.Range("H2").Value = 33.229931
$ws.Range("I2").Value = 0.06424093823669472
$ws.Range("J2").Value = 0.06424093823669472
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.490547666666667
$ws.Range("N2").Value = 19.471643
$ws.Range("O2").Value = 0.8021666724616637
$ws.Range("P2").Value = 0.8021666724616636
$ws.Range("Q2").Value = 71.89348370518145
$ws.Range("R2").Value = 647.0413533466331
$ws.Range("S2").Value = 0.05153193966114466
$ws.Range("T2").Value = 0.05153193966114465

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.07664366666667
$ws.Range("H3").Value = 33.229931
$ws.Range("I3").Value = 0.06424093823669472
$ws.Range("J3").Value = 0.06424093823669472
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.600723
$ws.Range("N3").Value = 4.802169
$ws.Range("O3").Value = 0.1978333275383364
$ws.Range("P3").Value = 0.1978333275383364
$ws.Range("Q3").Value = 17.73063828003767
$ws.Range("R3").Value = 159.575744520339
$ws.Range("S3").Value = 0.01270899857555007
$ws.Range("T3").Value = 0.01270899857555006

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 115.4636816666667
$ws.Range("H4").Value = 346.391045
$ws.Range("I4").Value = 0.6696518788314407
$ws.Range("J4").Value = 0.6696518788314407
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.490547666666667
$ws.Range("N4").Value = 19.471643
$ws.Range("O4").Value = 0.8021666724616637
$ws.Range("P4").Value = 0.8021666724616636
$ws.Range("Q4").Value = 749.4225296263261
$ws.Range("R4").Value = 6744.802766636934
$ws.Range("S4").Value = 0.537172419349918
$ws.Range("T4").Value = 0.537172419349918

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 115.4636816666667
$ws.Range("H5").Value = 346.391045
$ws.Range("I5").Value = 0.6696518788314407
$ws.Range("J5").Value = 0.6696518788314407
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.600723
$ws.Range("N5").Value = 4.802169
$ws.Range("O5").Value = 0.1978333275383364
$ws.Range("P5").Value = 0.1978333275383364
$ws.Range("Q5").Value = 184.8253709085117
$ws.Range("R5").Value = 1663.428338176605
$ws.Range("S5").Value = 0.1324794594815228
$ws.Range("T5").Value = 0.1324794594815228

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.82775733333333
$ws.Range("H6").Value = 113.483272
$ws.Range("I6").Value = 0.2193887151751843
$ws.Range("J6").Value = 0.2193887151751843
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.490547666666667
$ws.Range("N6").Value = 19.471643
$ws.Range("O6").Value = 0.8021666724616637
$ws.Range("P6").Value = 0.8021666724616636
$ws.Range("Q6").Value = 245.5228620950996
$ws.Range("R6").Value = 2209.705758855896
$ws.Range("S6").Value = 0.1759863156277173
$ws.Range("T6").Value = 0.1759863156277173

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 37.82775733333333
$ws.Range("H7").Value = 113.483272
$ws.Range("I7").Value = 0.2193887151751843
$ws.Range("J7").Value = 0.2193887151751843
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.600723
$ws.Range("N7").Value = 4.802169
$ws.Range("O7").Value = 0.1978333275383364
$ws.Range("P7").Value = 0.1978333275383364
$ws.Range("Q7").Value = 60.55176120188533
$ws.Range("R7").Value = 544.965850816968
$ws.Range("S7").Value = 0.04340239954746702
$ws.Range("T7").Value = 0.04340239954746702

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.055359000000001
$ws.Range("H8").Value = 24.166077
$ws.Range("I8").Value = 0.04671846775668023
$ws.Range("J8").Value = 0.04671846775668023
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.490547666666667
$ws.Range("N8").Value = 19.471643
$ws.Range("O8").Value = 0.8021666724616637
$ws.Range("P8").Value = 0.8021666724616636
$ws.Range("Q8").Value = 52.28369156161234
$ws.Range("R8").Value = 470.553224054511
$ws.Range("S8").Value = 0.03747599782288371
$ws.Range("T8").Value = 0.0374759978228837

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.055359000000001
$ws.Range("H9").Value = 24.166077
$ws.Range("I9").Value = 0.04671846775668023
$ws.Range("J9").Value = 0.04671846775668023
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.600723
$ws.Range("N9").Value = 4.802169
$ws.Range("O9").Value = 0.1978333275383364
$ws.Range("P9").Value = 0.1978333275383364
$ws.Range("Q9").Value = 12.894398424557
$ws.Range("R9").Value = 116.049585821013
$ws.Range("S9").Value = 0.009242469933796529
$ws.Range("T9").Value = 0.009242469933796528
